$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 'select host,user,AUTHENTICATION_STRING from mysql.user where user=''dcltest1'''
$ws.Range("L2").Value = 'select host,user,db,select_priv,insert_priv,update_priv,delete_priv from mysql.db where user=''dcltest1'' and host=''%'''
$ws.Range("N2").Value = 'select host,user,db,table_name,table_priv from mysql.tables_priv where user=''dcltest1'' and host=''%'''
$ws.Range("J3").Value = 'select host,user,AUTHENTICATION_STRING from mysql.user where user=''dcltest2'''
$ws.Range("L3").Value = 'select host,user,db,select_priv,insert_priv,update_priv,delete_priv from mysql.db where user=''dcltest2'' and host=''%'''
$ws.Range("N3").Value = 'select host,user,db,table_name,table_priv from mysql.tables_priv where user=''dcltest2'' and host=''%'''
$ws.Range("J4").Value = 'select host,user,AUTHENTICATION_STRING from mysql.user where user=''dcltest3'''
$ws.Range("L4").Value = 'select host,user,db,select_priv,insert_priv,update_priv,delete_priv from mysql.db where user=''dcltest3'' and host=''172.20.3.27'''
$ws.Range("N4").Value = 'select host,user,db,table_name,table_priv from mysql.tables_priv where user=''dcltest3'' and host=''172.20.3.27'''
$ws.Range("J5").Value = 'select user,AUTHENTICATION_STRING from mysql.user where user=''dcltest4'''
$ws.Range("L5").Value = 'select user,db,select_priv,insert_priv,update_priv,delete_priv from mysql.db where user=''dcltest4'''
$ws.Range("N5").Value = 'select user,db,table_name,table_priv from mysql.tables_priv where user=''dcltest4'''
$ws.Range("J6").Value = 'select host,user from mysql.user where user=''dcltest5'''
$ws.Range("L6").Value = 'select host,user,db,select_priv,insert_priv,update_priv,delete_priv from mysql.db where user=''dcltest5'' and host=''%'''
$ws.Range("N6").Value = 'select host,user,db,table_name,table_priv from mysql.tables_priv where user=''dcltest5'' and host=''%'''
$ws.Range("J7").Value = 'select host,user,select_priv,insert_priv,update_priv,delete_priv,create_priv,drop_priv from mysql.user where user=''dcltest6'' and host=''%'''
$ws.Range("L7").Value = 'select host,user,db,select_priv,insert_priv,update_priv,delete_priv,create_priv,drop_priv from mysql.db where user=''dcltest6'' and host=''%'''
$ws.Range("N7").Value = 'select host,user,db,table_name,table_priv from mysql.tables_priv where user=''dcltest6'' and host=''%'''
$ws.Range("J8").Value = 'select host,user,select_priv,insert_priv,update_priv,delete_priv,create_priv,drop_priv,reload_priv,shutdown_priv,PROCESS_PRIV,FILE_PRIV,GRANT_PRIV,REFERENCES_PRIV,INDEX_PRIV,SHOW_DB_PRIV,CREATE_USER_PRIV from mysql.user where user=''dcltest7'''
$ws.Range("L8").Value = 'select host,user,db,select_priv,insert_priv,update_priv,delete_priv,create_priv,drop_priv from mysql.db where user=''dcltest7'' and host=''%'''
$ws.Range("N8").Value = 'select host,user,db,table_name,table_priv from mysql.tables_priv where user=''dcltest7'' and host=''%'''
$ws.Range("J10").Value = 'select host,user from mysql.user where user=''dcltest9'''
$ws.Range("L10").Value = 'select host,user,db,select_priv,insert_priv,update_priv,delete_priv,create_priv,drop_priv from mysql.db where user=''dcltest9'' and host=''%'''
$ws.Range("N10").Value = 'select host,user,db,table_name,table_priv from mysql.tables_priv where user=''dcltest9'' and host=''%'''

$ws.Range("J18").Select()
